$d = $word.ActiveDocument

# 1. Insert a new bullet paragraph before the first paragraph
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs(1)
$newPara.Range.Text = "the price doesn't calculate the sub element prices and when they're unchecked it zeros all togehter"

# 2. Merge the split runs "Add main heading class for centered h" + "2 " + "with increased font-sized"
#    into a single run by finding the full text and replacing it with itself.
$text2 = "Add main heading class for centered h2 with increased font-sized"
$d.Content.Find.Execute($text2, $false, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# 3. Merge the split runs "Calculate price is not accurate it still doesn" + "’" + "t include quantity of items"
$text3 = "Calculate price is not accurate it still doesn’t include quantity of items"
$d.Content.Find.Execute($text3, $false, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

# 4. Merge the split runs "Ingredient box is not well designed and not in all " + "3 " + "items in special for the week section"
$text4 = "Ingredient box is not well designed and not in all 3 items in special for the week section"
$d.Content.Find.Execute($text4, $false, $false, $false, $false, $false, $true, 1, $false, $text4, 2) | Out-Null

Write-Host "done"
